$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - first sale entry
$ws.Range("A2").Value = "14-07-2021"
$ws.Range("B2").Value = 0
$ws.Range("C2").Value = "Item A"
$ws.Range("D2").Value = 1
$ws.Range("E2").Value = 10
$ws.Range("F2").Value = 10

# Row 3 - second sale entry
$ws.Range("A3").Value = "14-07-2021"
$ws.Range("B3").Value = 2
$ws.Range("C3").Value = "Item C"
$ws.Range("D3").Value = 2
$ws.Range("E3").Value = 10
$ws.Range("F3").Value = 20
